# Auto-committed edit: update DB-layout remark text for PfReward.xlsx
# - RepayType remark (G14): reformat "0.撥款... 2.部分償還 3.提前結案" to
#   colon-separated, one item per line.
# - PerfDate remark (G10): change "1.xxx" / "2.xxx" style numbering to
#   "1:xxx" / "2:xxx".
# - IntroducerBonus / IntroducerAddBonus / CoorgnizerBonus remarks (G21,
#   G23, G25): change half-width colon "薪碼:Qx" to full-width "薪碼：Qx".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("G10").Value = "1:撥貸(計件代碼變更)，為撥款日期`n2:部分償還、提前結案，為會計日"
$ws.Range("G14").Value = "0:撥款(計件代碼變更)`n2:部分償還`n3:提前結案"
$ws.Rows.Item(14).RowHeight = 48.6
$ws.Range("G21").Value = "薪碼：Q2"
$ws.Range("G23").Value = "薪碼：Q1"
$ws.Range("G25").Value = "薪碼：Q2"

$ws.Range("G26").Select()

Write-Output "PfReward.xlsx edits applied"
